$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2125
$ws.Range("J97").Value = 2125
$ws.Range("L97").Value = 6375
$ws.Range("N97").Value = -7367
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H137").Value = 16282476
$ws.Range("I137").Value = 3907127.5
$ws.Range("K137").Value = 11721382.5
$ws.Range("M137").Value = -11718832.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14349.423
$ws.Range("I2").Value = 20228.389
$ws.Range("J2").Value = 1121.75
$ws.Range("K2").Value = 20228.389
$ws.Range("L2").Value = 1121.75
$ws.Range("M2").Value = -20115.389
$ws.Range("N2").Value = -1347.75
$ws.Range("H61").Value = 6281637.5
$ws.Range("I61").Value = 2778845
$ws.Range("J61").Value = 58823530
$ws.Range("K61").Value = 2778845
$ws.Range("L61").Value = 58823530
$ws.Range("M61").Value = -2778633
$ws.Range("N61").Value = -58823954
$ws.Range("H101").Value = 35500
$ws.Range("J101").Value = 35500
$ws.Range("L101").Value = 35500
$ws.Range("N101").Value = -41990
$ws.Range("H116").Value = 14349.423
$ws.Range("I116").Value = 20228.389
$ws.Range("J116").Value = 1121.75
$ws.Range("K116").Value = 20228.389
$ws.Range("L116").Value = 1121.75
$ws.Range("M116").Value = -17934.389
$ws.Range("N116").Value = -5709.75
$ws.Range("H132").Value = 6947116
$ws.Range("I132").Value = 6805437
$ws.Range("J132").Value = 7938867
$ws.Range("K132").Value = 20416311
$ws.Range("L132").Value = 23816601
$ws.Range("M132").Value = -20413781
$ws.Range("N132").Value = -23821661
$ws.Range("H136").Value = 6281637.5
$ws.Range("I136").Value = 2778845
$ws.Range("J136").Value = 58823530
$ws.Range("K136").Value = 8336535
$ws.Range("L136").Value = 176470590
$ws.Range("M136").Value = -8333985
$ws.Range("N136").Value = -176475690

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14349.423
$ws.Range("I3").Value = 20228.389
$ws.Range("J3").Value = 1121.75
$ws.Range("K3").Value = 20228.389
$ws.Range("L3").Value = 1121.75
$ws.Range("M3").Value = -20114.389
$ws.Range("N3").Value = -1349.75
$ws.Range("H22").Value = 503.17648
$ws.Range("I22").Value = 503.17648
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 503.17648
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -330.17648
$ws.Range("N22").ClearContents()
$ws.Range("H105").Value = 1966.4117
$ws.Range("I105").Value = 1939.091
$ws.Range("J105").Value = 2016.5
$ws.Range("K105").Value = 1939.091
$ws.Range("L105").Value = 2016.5
$ws.Range("M105").Value = -192.0909999999999
$ws.Range("N105").Value = -5510.5
$ws.Range("H107").Value = 933.3333
$ws.Range("I107").Value = 950
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 950
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 970
$ws.Range("N107").Value = -4740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3666920.5
$ws.Range("I70").Value = 1857116.8
$ws.Range("J70").Value = 9096332
$ws.Range("K70").Value = 1857116.8
$ws.Range("L70").Value = 9096332
$ws.Range("M70").Value = -1856846.8
$ws.Range("N70").Value = -9096872
$ws.Range("H73").Value = 3666920.5
$ws.Range("I73").Value = 1857116.8
$ws.Range("J73").Value = 9096332
$ws.Range("K73").Value = 1857116.8
$ws.Range("L73").Value = 9096332
$ws.Range("M73").Value = -1856180.8
$ws.Range("N73").Value = -9098204
$ws.Range("H97").Value = 15626511
$ws.Range("I97").Value = 965.3077
$ws.Range("K97").Value = 965.3077
$ws.Range("M97").Value = -469.3077
$ws.Range("H113").Value = 29984
$ws.Range("I113").Value = 6170
$ws.Range("J113").Value = 53798
$ws.Range("K113").Value = 6170
$ws.Range("L113").Value = 53798
$ws.Range("M113").Value = -4000
$ws.Range("N113").Value = -58138
$ws.Range("H122").Value = 33436576
$ws.Range("I122").Value = 170402.33
$ws.Range("J122").Value = 83335840
$ws.Range("K122").Value = 511206.99
$ws.Range("L122").Value = 250007520
$ws.Range("M122").Value = -508756.99
$ws.Range("N122").Value = -250012420

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1150.6316
$ws.Range("I7").Value = 1044.1333
$ws.Range("J7").Value = 1550
$ws.Range("K7").Value = 1044.1333
$ws.Range("L7").Value = 1550
$ws.Range("M7").Value = -932.1333
$ws.Range("N7").Value = -1774
$ws.Range("H22").Value = 5428.4
$ws.Range("I22").Value = 4760
$ws.Range("J22").Value = 5762.6
$ws.Range("K22").Value = 4760
$ws.Range("L22").Value = 5762.6
$ws.Range("M22").Value = -4465
$ws.Range("N22").Value = -6352.6
$ws.Range("H27").Value = 5428.4
$ws.Range("I27").Value = 4760
$ws.Range("J27").Value = 5762.6
$ws.Range("K27").Value = 4760
$ws.Range("L27").Value = 5762.6
$ws.Range("M27").Value = -4653
$ws.Range("N27").Value = -5976.6
$ws.Range("H61").Value = 1099
$ws.Range("I61").Value = 970
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 970
$ws.Range("L61").Value = 1400
$ws.Range("M61").Value = -768
$ws.Range("N61").Value = -1804
$ws.Range("H100").Value = 1721.3334
$ws.Range("I100").Value = 1299.5555
$ws.Range("J100").Value = 2986.6667
$ws.Range("K100").Value = 1299.5555
$ws.Range("L100").Value = 2986.6667
$ws.Range("M100").Value = -758.5554999999999
$ws.Range("N100").Value = -4068.6667
$ws.Range("H113").Value = 1099
$ws.Range("I113").Value = 970
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 970
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 1200
$ws.Range("N113").Value = -5740
$ws.Range("H126").Value = 1150.6316
$ws.Range("I126").Value = 1044.1333
$ws.Range("J126").Value = 1550
$ws.Range("K126").Value = 3132.3999
$ws.Range("L126").Value = 4650
$ws.Range("M126").Value = -662.3998999999999
$ws.Range("N126").Value = -9590
